$wb = $excel.ActiveWorkbook

# Sheet "Final_Total" updates (2023 tax year rows)
$wsTotal = $wb.Worksheets.Item("Final_Total")
$wsTotal.Range("D47").Value = 13590.8188469453
$wsTotal.Range("D48").Value = 56842.5511530548
$wsTotal.Range("D49").Value = 402.716224242701
$wsTotal.Range("D50").Value = 2728.61384526203
$wsTotal.Range("D51").Value = 746.939930495271

# Sheet "Final_Gahanna" updates (2023 tax year rows)
$wsGahanna = $wb.Worksheets.Item("Final_Gahanna")
$wsGahanna.Range("D20").Value = 13590.8188469453
$wsGahanna.Range("D21").Value = 56842.5511530548

# Sheet "Final_JeffersonUnincorporated" updates (2023 tax year rows)
$wsJefferson = $wb.Worksheets.Item("Final_JeffersonUnincorporated")
$wsJefferson.Range("D29").Value = 402.716224242701
$wsJefferson.Range("D30").Value = 2728.61384526203
$wsJefferson.Range("D31").Value = 746.939930495271
